$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width (OOXML width 14.7109375 -> 16.42578125)
$ws.Columns.Item(2).ColumnWidth = 15.6

# Update cell values in columns A and B for rows 1-32
$ws.Range("A1").Value = -0.3632265282951721
$ws.Range("B1").Value = 0.36242495512692585
$ws.Range("A2").Value = -0.23833405228415927
$ws.Range("B2").Value = 0.23615228160239887
$ws.Range("A3").Value = -0.13320302085143965
$ws.Range("B3").Value = 0.13268541695299874
$ws.Range("A4").Value = -0.12068541701465918
$ws.Range("B4").Value = 0.12024544079369548
$ws.Range("A5").Value = -0.11424544100713518
$ws.Range("B5").Value = 0.11339050250879801
$ws.Range("A6").Value = -0.039447609305863551
$ws.Range("B6").Value = 0.039425321807485147
$ws.Range("A7").Value = -0.019425322071766615
$ws.Range("B7").Value = 0.019406550562090175
$ws.Range("A8").Value = 0.00059344917301284994
$ws.Range("B8").Value = -0.00059459605551115402
$ws.Range("A9").Value = 0.006594595836177497
$ws.Range("B9").Value = -0.0066001889007987202
$ws.Range("A10").Value = 0.012600188681801683
$ws.Range("B10").Value = -0.012600500909989876
$ws.Range("A11").Value = 0.017100500695988785
$ws.Range("B11").Value = -0.017110357414289012
$ws.Range("A12").Value = 0.023110357195469611
$ws.Range("B12").Value = -0.023183255086506538
$ws.Range("A13").Value = -0.039147876104077106
$ws.Range("B13").Value = 0.039083142654182446
$ws.Range("A14").Value = -0.027083142891765277
$ws.Range("B14").Value = 0.02705171589399491
$ws.Range("A15").Value = -0.021051716112845398
$ws.Range("B15").Value = 0.021026962333941412
$ws.Range("A16").Value = -0.015026962553537526
$ws.Range("B16").Value = 0.015004479209294264
$ws.Range("A17").Value = -0.0090044794298407282
$ws.Range("B17").Value = 0.0089999997694336642
$ws.Range("A18").Value = -0.036111730316036983
$ws.Range("B18").Value = 0.036097206721326103
$ws.Range("A19").Value = -0.027097206934413531
$ws.Range("B19").Value = 0.027014036881643833
$ws.Range("A20").Value = -0.018014037096673263
$ws.Range("B20").Value = 0.018004328295976535
$ws.Range("A21").Value = -0.0090043285113203808
$ws.Range("B21").Value = 0.0089999997844190105
$ws.Range("A22").Value = -0.093949260969612425
$ws.Range("B22").Value = 0.093635436078372436
$ws.Range("A23").Value = -0.084635436299084432
$ws.Range("B23").Value = 0.084127100331170546
$ws.Range("A24").Value = -0.042127100662946226
$ws.Range("B24").Value = 0.041999999666370691
$ws.Range("A25").Value = -0.067942893421328421
$ws.Range("B25").Value = 0.067865330144218916
$ws.Range("A26").Value = -0.061865330363175985
$ws.Range("B26").Value = 0.061773072114728933
$ws.Range("A27").Value = -0.055773072334611484
$ws.Range("B27").Value = 0.055484609888774461
$ws.Range("A28").Value = -0.075266026792873753
$ws.Range("B28").Value = 0.074533089593895596
$ws.Range("A29").Value = -0.062533089839826417
$ws.Range("B29").Value = 0.062169284815002612
$ws.Range("A30").Value = -0.042169285088635267
$ws.Range("B30").Value = 0.042020275197652612
$ws.Range("A31").Value = -0.027020275457367404
$ws.Range("B31").Value = 0.02700095312320272
$ws.Range("A32").Value = -0.0060009534028262479
$ws.Range("B32").Value = 0.0059999997690614038
